# changes for minimum threshold
# Adds a "Minimum_work" threshold column and a "Feedback" column to the
# SheetGradingOrder sheet, marking each gallery-key sheet with a minimum
# work threshold of 10 and a "Need more work" feedback message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetGradingOrder")

# Header row
$ws.Range("C1").Value = "Minimum_work"
$ws.Range("D1").Value = "Feedback"

# Constant Samples row
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = "Need more work"

# Formula Samples row
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "Need more work"

# SoftFormula Samples row
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = "Need more work"

# Test Case Samples row (no feedback message for this one)
$ws.Range("C5").Value = 10

$ws.Range("C6").Select()
